$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price-record row was inserted above the existing row 433,
# pushing the old rows 433:518 down to 434:519 (dimension A1:R518 -> A1:R519).
$ws.Rows("433").Insert()

$ws.Range("A433").Value = 10
$ws.Range("B433").Value = "Vega Modelo de Temuco"
$ws.Range("C433").Value = "La Araucanía"
$ws.Range("D433").Value = 45209
$ws.Range("E433").Value = 9
$ws.Range("F433").Value = 100112001
$ws.Range("G433").Value = "Berenjena"
$ws.Range("H433").Value = "Sin especificar"
$ws.Range("I433").Value = "Primera"
$ws.Range("J433").Value = 100
$ws.Range("K433").Value = 12000
$ws.Range("L433").Value = 12000
$ws.Range("M433").Value = 12000
$ws.Range("N433").Value = "`$/caja 40 unidades"
$ws.Range("O433").Value = "Región de Arica y Parinacota"
$ws.Range("P433").Value = 300
$ws.Range("Q433").Value = 40
$ws.Range("R433").Value = "Hortaliza"
